$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename sheet ---
$ws.Name = "Datos_Sujeto3_Dummies"

# --- Header row updates (E1:H1 relabeled, I1/J1 new) ---
$ws.Range("E1").Value = "Corto_enCorto"
$ws.Range("F1").Value = "Corto_enLargo"
$ws.Range("G1").Value = "Largo_enCorto"
$ws.Range("H1").Value = "Largo_enLargo"
$ws.Range("I1").Value = "EnsayosCortos"
$ws.Range("J1").Value = "EnsayosLargos"

# --- Column widths for new/relabeled columns (auto-fit-like sizing) ---
$ws.Columns.Item(5).ColumnWidth = 13.5834
$ws.Columns.Item(6).ColumnWidth = 12.5834
$ws.Columns.Item(7).ColumnWidth = 12.5834
$ws.Columns.Item(8).ColumnWidth = 12.75
$ws.Columns.Item(9).ColumnWidth = 13.0834
$ws.Columns.Item(10).ColumnWidth = 12.0834

# --- Shift data: old F (Rechazos), G (FA), H (Signals) columns ---
# New layout: F <- old G, G <- (cleared), H <- old F, I <- old H
$lastRow = 81
for ($r = 2; $r -le $lastRow; $r++) {
    $oldF = $ws.Cells.Item($r, 6).Value2
    $oldG = $ws.Cells.Item($r, 7).Value2
    $oldH = $ws.Cells.Item($r, 8).Value2

    $ws.Cells.Item($r, 6).Value = $oldG
    $ws.Cells.Item($r, 7).ClearContents()
    $ws.Cells.Item($r, 8).Value = $oldF
    $ws.Cells.Item($r, 9).Value = $oldH
}

# --- New J column: EnsayosLargos = F + H, filled down ---
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 10).Formula = "=F$r+H$r"
}

# --- Selection matches post-edit state ---
$ws.Range("J2:J81").Select() | Out-Null
